$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that needs to move
# forward by one day (45204 -> 45205) for every data row (rows 2-540).
$ws.Range("C2:C540").Value = 45205
